# Commit #5: property boat&car done
# Fix the "汽車" (car) sheet (Worksheets item 3): row 1 had accidentally been
# filled with a copy of the data row instead of real column headers, and the
# sheet was missing the trailing property_category/category/date/
# legislator_name/legislator_id/source_file/index columns that every other
# property sheet already carries. Rebuild row 1 as headers and extend both
# rows out to column N.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Row 1: replace the bogus duplicated data with real column headers ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"

# New header cells H1:N1 - clone formatting from the existing bold/boxed
# header style (B1) before writing the label text.
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2: keep the existing A2:G2 data, append the missing columns ---
# New data cells H2:N2 - clone formatting from an existing data cell (B2)
# before writing the values.
$ws.Range("B2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)

$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2012-04-30"
$ws.Range("K2").Value = "楊麗環"
$ws.Range("L2").Value = 960
$ws.Range("M2").Value = "tmp700a1"
$ws.Range("N2").Value = 30
